$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.242.79'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.829.60'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  +0.50%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.69'
$ws.Range('E5').Value = '  -1.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5985'
$ws.Range('E6').Value = '  -4.46%  '
$ws.Range('E7').Value = '  +0.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.06973'
$ws.Range('E8').Value = '  -5.85%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2758'
$ws.Range('E9').Value = '  -4.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.31'
$ws.Range('E10').Value = '  -6.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07634'
$ws.Range('E11').Value = '  -1.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.835.88'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.775'
$ws.Range('E13').Value = '  -3.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6263'
$ws.Range('E14').Value = '  -7.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009710'
$ws.Range('E15').Value = '  -5.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '78.61'
$ws.Range('E16').Value = '  -3.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.056.76'
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.748'
$ws.Range('E18').Value = '  -7.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '221.81'
$ws.Range('E19').Value = '  -5.32%  '
$ws.Range('E20').Value = '  +0.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.58'
$ws.Range('E21').Value = '  -6.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.907'
$ws.Range('E22').Value = '  -5.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.006'
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '156.06'
$ws.Range('E24').Value = '  -1.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1296'
$ws.Range('E25').Value = '  -3.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.955'
$ws.Range('E26').Value = '  -6.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.52'
$ws.Range('E27').Value = '  -4.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06768'
$ws.Range('E28').Value = '  -6.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.447'
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.443'
$ws.Range('E30').Value = '  -2.40%  '
$ws.Range('E31').Value = '  -4.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.769'
$ws.Range('E32').Value = '  -7.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.093'
$ws.Range('E33').Value = '  -4.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.723'
$ws.Range('E34').Value = '  -5.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6441'
$ws.Range('E35').Value = '  -7.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.546'
$ws.Range('E36').Value = '  -0.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.738'
$ws.Range('E37').Value = '  -2.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.195.74'
$ws.Range('E38').Value = '  -2.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01738'
$ws.Range('E39').Value = '  -5.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.509'
$ws.Range('E40').Value = '  -5.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9026'
$ws.Range('E41').Value = '  -5.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.004'
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.983.61'
$ws.Range('E43').Value = '  -0.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.33'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.37'
$ws.Range('E45').Value = '  -4.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000113'
$ws.Range('E46').Value = '  -2.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.507'
$ws.Range('E47').Value = '  -3.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.576'
$ws.Range('E48').Value = '  -7.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4556'
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05515'
$ws.Range('E50').Value = '  -2.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.374'
$ws.Range('E51').Value = '  -8.12%  '
